$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.746.06"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.924.57"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.61"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.69"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.571"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.61"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0847"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "3.382.48"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.99"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.05"
$ws.Range("E16").Value = "  +65.79%  "
$ws.Range("D17").Value = "2.933.70"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.989"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "50.704.53"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.00"
$ws.Range("E20").Value = "  -6.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "0.0₃0943"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.33"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.60"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  +11.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.84"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.36"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.95"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.39"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.11"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0429"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.28"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "123.02"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.91"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D48").Value = "1.998.67"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.257"
$ws.Range("E49").Value = "  -5.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0314"
$ws.Range("E50").Value = "  -4.75%  "
$ws.Range("E51").Value = "  +4.17%  "
